$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$cb = $wb.Worksheets.Item("Codebook")

# ------------------------------------------------------------------
# The "Data" sheet columns C:E got rearranged (Gender moved before Age,
# the household-size column was renamed to NPH and moved to column D,
# and several rows of Gender/NPH data were corrected). Apply the new
# formatting first (copied from stable, untouched reference cells),
# then write the corrected values.
# ------------------------------------------------------------------

# Stable style reference cells that never change: A1 (bold header),
# A2 (plain data), E2 (centred data) -- style-copy from these BEFORE
# any values are overwritten.
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

$ws.Range("E2").Copy()
$ws.Range("D1:D15").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$ws.Range("C2:C15").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# Corrected values for columns C, D, E (A and B are untouched)
# ------------------------------------------------------------------
$ws.Range("C1").Value = "Gender"
$ws.Range("D1").Value = "NPH"
$ws.Range("E1").Value = "Age"

$ws.Range("C2").Value = "M"
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 23

$ws.Range("C3").Value = "O"
$ws.Range("D3").Value = 9
$ws.Range("E3").Value = 34

$ws.Range("C4").Value = "F"
$ws.Range("D4").Value = "NA"
$ws.Range("E4").Value = 64

$ws.Range("C5").Value = "F"
$ws.Range("D5").Value = 8
$ws.Range("E5").Value = 41

$ws.Range("C6").Value = "NA"
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = 20

$ws.Range("C7").Value = "F"
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = 77

$ws.Range("C8").Value = "O"
$ws.Range("D8").Value = 7
$ws.Range("E8").Value = 0

$ws.Range("C9").Value = "M"
$ws.Range("D9").Value = 7
$ws.Range("E9").Value = 75

$ws.Range("C10").Value = "N"
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 24

$ws.Range("C11").Value = "M"
$ws.Range("D11").Value = 5
$ws.Range("E11").Value = 33

$ws.Range("C12").Value = "F"
$ws.Range("D12").Value = 9
$ws.Range("E12").Value = 59

$ws.Range("C13").Value = "F"
$ws.Range("D13").Value = 5
$ws.Range("E13").Value = 75

$ws.Range("C14").Value = "M"
$ws.Range("D14").Value = 10
$ws.Range("E14").Value = 0

$ws.Range("C15").Value = "M"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 72

# ------------------------------------------------------------------
# Codebook sheet: the household-size row now references the new NPH name
# ------------------------------------------------------------------
$cb.Range("A3").Value = "Number of people in household(NPH)"

# ------------------------------------------------------------------
# Restore the last-used selection on the Data sheet
# ------------------------------------------------------------------
$ws.Activate()
$ws.Range("L10").Select() | Out-Null
